# Auto-generated: apply cryptos.xlsx crypto-price refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'59.240.11"
$ws.Range('E2').Value = '  +0.81%  '

$ws.Range('D3').Value = "'2.602.08"
$ws.Range('E3').Value = '  +0.58%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').Value = "'540.65"
$ws.Range('E5').Value = '  +4.03%  '

$ws.Range('D6').Value = "'141.19"
$ws.Range('E6').Value = '  +1.05%  '

$ws.Range('E7').Value = '  +0.10%  '

$ws.Range('E8').Value = '  +0.38%  '

$ws.Range('D9').Value = "'6.45"
$ws.Range('E9').Value = '  -0.59%  '

$ws.Range('E10').Value = '  +2.18%  '

$ws.Range('E11').Value = '  +1.19%  '

$ws.Range('E12').Value = '  +1.22%  '

$ws.Range('D13').Value = "'3.062.76"
$ws.Range('E13').Value = '  +0.61%  '

$ws.Range('D14').Value = "'59.157.87"
$ws.Range('E14').Value = '  +0.58%  '

$ws.Range('D15').Value = "'20.57"
$ws.Range('E15').Value = '  +0.56%  '

$ws.Range('D16').Value = "'2.651.49"
$ws.Range('E16').Value = '  +2.64%  '

$ws.Range('E17').Value = '  +0.68%  '

$ws.Range('D18').Value = "'341.57"
$ws.Range('E18').Value = '  +0.74%  '

$ws.Range('D19').Value = "'4.37"
$ws.Range('E19').Value = '  +1.00%  '

$ws.Range('D20').Value = "'10.12"
$ws.Range('E20').Value = '  -0.76%  '

$ws.Range('E21').Value = '  -2.01%  '

$ws.Range('D22').Value = "'1.00"
$ws.Range('E22').Value = '  +0.09%  '

$ws.Range('D23').Value = "'67.63"
$ws.Range('E23').Value = '  +2.14%  '

$ws.Range('D24').Value = "'0.408"
$ws.Range('E24').Value = '  +0.88%  '

$ws.Range('E25').Value = '  -0.76%  '

$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = '  +0.19%  '

$ws.Range('E27').Value = '  +1.60%  '

$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').Value = "'0.0₃0752"
$ws.Range('E28').Value = '  +4.02%  '

$ws.Range('B29').Value = 'USDe'
$ws.Range('C29').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D29').Value = "'0.999"
$ws.Range('E29').Value = '  +0.04%  '

$ws.Range('E30').Value = '  +8.96%  '

$ws.Range('D31').Value = "'5.81"
$ws.Range('E31').Value = '  -2.28%  '

$ws.Range('D32').Value = "'18.72"
$ws.Range('E32').Value = '  -0.31%  '

$ws.Range('D33').Value = "'149.43"
$ws.Range('E33').Value = '  +0.48%  '

$ws.Range('D34').Value = "'3.99"
$ws.Range('E34').Value = '  +0.04%  '

$ws.Range('D35').Value = "'37.14"
$ws.Range('E35').Value = '  +2.41%  '

$ws.Range('E36').Value = '  -1.27%  '

$ws.Range('E37').Value = '  +0.77%  '

$ws.Range('D38').Value = "'0.836"
$ws.Range('E38').Value = '  +0.20%  '

$ws.Range('D39').Value = "'0.818"
$ws.Range('E39').Value = '  +0.03%  '

$ws.Range('E40').Value = '  +1.55%  '

$ws.Range('D41').Value = "'0.999"
$ws.Range('E41').Value = '  +0.16%  '

$ws.Range('D42').Value = "'274.03"
$ws.Range('E42').Value = '  -0.32%  '

$ws.Range('E43').Value = '  +1.63%  '

$ws.Range('D44').Value = "'10.74"
$ws.Range('E44').Value = '  -0.14%  '

$ws.Range('D45').Value = "'0.0957"
$ws.Range('E45').Value = '  +0.66%  '

$ws.Range('D46').Value = "'0.0525"
$ws.Range('E46').Value = '  +0.66%  '

$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = "'1.944.75"
$ws.Range('E47').Value = '  -1.90%  '

$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = "'0.0223"
$ws.Range('E48').Value = '  +1.20%  '

$ws.Range('D49').Value = "'18.48"
$ws.Range('E49').Value = '  +3.07%  '

$ws.Range('E50').Value = '  +0.64%  '

$ws.Range('D51').Value = "'111.56"
$ws.Range('E51').Value = '  -1.13%  '
